$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest cryptos snapshot refreshed by the scheduled GitHub Actions job.
# Column D holds price text (e.g. "27.154.97"); a leading apostrophe forces
# Excel to keep it as text instead of silently parsing it as a Double.

$ws.Range("D2").Value = "'" + '27.154.97'
$ws.Range("E2").Value = '  +0.59%  '

$ws.Range("D3").Value = "'" + '1.901.67'
$ws.Range("E3").Value = '  +1.05%  '

$ws.Range("D4").Value = "'" + '1.001'
$ws.Range("E4").Value = '  +0.37%  '

$ws.Range("D5").Value = "'" + '305.91'
$ws.Range("E5").Value = '  -0.20%  '

$ws.Range("D6").Value = "'" + '1.001'
$ws.Range("E6").Value = '  +0.32%  '

$ws.Range("D7").Value = "'" + '0.5258'
$ws.Range("E7").Value = '  +1.62%  '

$ws.Range("D8").Value = "'" + '0.3774'
$ws.Range("E8").Value = '  +1.60%  '

$ws.Range("D9").Value = "'" + '0.07250'
$ws.Range("E9").Value = '  +0.80%  '

$ws.Range("D10").Value = "'" + '21.18'
$ws.Range("E10").Value = '  +1.61%  '

$ws.Range("D11").Value = "'" + '0.8982'
$ws.Range("E11").Value = '  -0.36%  '

$ws.Range("D12").Value = "'" + '0.08333'
$ws.Range("E12").Value = '  +10.24%  '

$ws.Range("D13").Value = "'" + '1.891.47'
$ws.Range("E13").Value = '  -2.14%  '

$ws.Range("D14").Value = "'" + '94.83'
$ws.Range("E14").Value = '  -0.48%  '

$ws.Range("D15").Value = "'" + '5.268'
$ws.Range("E15").Value = '  +0.22%  '

$ws.Range("D16").Value = "'" + '1.001'
$ws.Range("E16").Value = '  +0.38%  '

$ws.Range("D17").Value = "'" + '0.000008614'
$ws.Range("E17").Value = '  +1.34%  '

$ws.Range("D18").Value = "'" + '14.50'
$ws.Range("E18").Value = '  +1.77%  '

$ws.Range("D19").Value = "'" + '1.000'
$ws.Range("E19").Value = '  +0.23%  '

$ws.Range("D20").Value = "'" + '27.190.80'
$ws.Range("E20").Value = '  +0.55%  '

$ws.Range("D21").Value = "'" + '5.060'
$ws.Range("E21").Value = '  +0.60%  '

$ws.Range("D22").Value = "'" + '2.135.82'
$ws.Range("E22").Value = '  +1.04%  '

$ws.Range("D23").Value = "'" + '10.58'
$ws.Range("E23").Value = '  +1.63%  '

$ws.Range("D24").Value = "'" + '6.430'
$ws.Range("E24").Value = '  -0.33%  '

$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = "'" + '146.46'
$ws.Range("E25").Value = '  +0.55%  '

$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = "'" + '2.278'
$ws.Range("E26").Value = '  +7.31%  '

$ws.Range("D27").Value = "'" + '1.758'
$ws.Range("E27").Value = '  -1.57%  '

$ws.Range("D28").Value = "'" + '18.10'
$ws.Range("E28").Value = '  +0.57%  '

$ws.Range("D29").Value = "'" + '114.77'
$ws.Range("E29").Value = '  +0.20%  '

$ws.Range("D30").Value = "'" + '4.927'
$ws.Range("E30").Value = '  -0.35%  '

$ws.Range("D31").Value = "'" + '4.779'
$ws.Range("E31").Value = '  -0.08%  '

$ws.Range("D32").Value = "'" + '0.09256'
$ws.Range("E32").Value = '  +0.61%  '

$ws.Range("D33").Value = "'" + '0.8135'
$ws.Range("E33").Value = '  +6.89%  '

$ws.Range("D34").Value = "'" + '0.05049'
$ws.Range("E34").Value = '  +0.34%  '

$ws.Range("D35").Value = "'" + '1.235'
$ws.Range("E35").Value = '  +4.08%  '

$ws.Range("D36").Value = "'" + '2.981'
$ws.Range("E36").Value = '  -0.66%  '

$ws.Range("D37").Value = "'" + '3.332'
$ws.Range("E37").Value = '  +1.56%  '

$ws.Range("D38").Value = "'" + '2.583'
$ws.Range("E38").Value = '  +3.09%  '

$ws.Range("D39").Value = "'" + '0.5706'
$ws.Range("E39").Value = '  +2.12%  '

$ws.Range("D40").Value = "'" + '0.01978'
$ws.Range("E40").Value = '  -0.66%  '

$ws.Range("E41").Value = '  -0.07%  '

$ws.Range("D42").Value = "'" + '6.661'
$ws.Range("E42").Value = '  +1.15%  '

$ws.Range("D43").Value = "'" + '8.948'
$ws.Range("E43").Value = '  +1.14%  '

$ws.Range("E44").Value = '  +0.62%  '

$ws.Range("D45").Value = "'" + '0.1511'
$ws.Range("E45").Value = '  +0.55%  '

$ws.Range("D46").Value = "'" + '0.4836'
$ws.Range("E46").Value = '  +1.10%  '

$ws.Range("D47").Value = "'" + '1.001'
$ws.Range("E47").Value = '  +0.39%  '

$ws.Range("D48").Value = "'" + '10.13'
$ws.Range("E48").Value = '  -0.19%  '

$ws.Range("D49").Value = "'" + '1.614'
$ws.Range("E49").Value = '  +2.97%  '

$ws.Range("D50").Value = "'" + '37.45'
$ws.Range("E50").Value = '  +1.03%  '

$ws.Range("D51").Value = "'" + '63.51'
$ws.Range("E51").Value = '  +0.16%  '
